# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型" sheets
# to reflect refreshed stats from the bilibili source, per commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 184
$wsExhibit.Range("F4").Value = 340
$wsExhibit.Range("F5").Value = 404
$wsExhibit.Range("F7").Value = 2376
$wsExhibit.Range("F8").Value = 401
$wsExhibit.Range("F9").Value = 6080

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 184
$wsAll.Range("F4").Value = 340
$wsAll.Range("F5").Value = 404
$wsAll.Range("F9").Value = 2376
$wsAll.Range("F10").Value = 401
$wsAll.Range("F11").Value = 6080
